# Auto-generated edit script applying scheduled Kraken_Profits value refresh
# across the ALC/ARM/BSM/CUL/GSM/LTW/WVR sheets (CRP untouched).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 333.33334
$ws.Range("H80").Value = 5735
$ws.Range("I80").Value = 633.3333
$ws.Range("J80").Value = 13387.5
$ws.Range("K80").Value = 1899.9999
$ws.Range("L80").Value = 40162.5
$ws.Range("M80").Value = -901.9999
$ws.Range("N80").Value = -42158.5
$ws.Range("H83").Value = 5735
$ws.Range("I83").Value = 633.3333
$ws.Range("J83").Value = 13387.5
$ws.Range("K83").Value = 5699.9997
$ws.Range("L83").Value = 120487.5
$ws.Range("M83").Value = -707.9997000000003
$ws.Range("N83").Value = -130471.5
$ws.Range("H96").Value = 3416.5
$ws.Range("I96").Value = 3744.8
$ws.Range("K96").Value = 11234.4
$ws.Range("M96").Value = -9861.400000000001
$ws.Range("H98").Value = 1708.091
$ws.Range("I98").Value = 1433.9
$ws.Range("J98").Value = 4450
$ws.Range("K98").Value = 1433.9
$ws.Range("L98").Value = 4450
$ws.Range("M98").Value = 64.09999999999991
$ws.Range("N98").Value = -7446
$ws.Range("H111").Value = 5216.6665
$ws.Range("I111").Value = 4460
$ws.Range("J111").Value = 9000
$ws.Range("K111").Value = 13380
$ws.Range("L111").Value = 27000
$ws.Range("M111").Value = -10313
$ws.Range("N111").Value = -33134
$ws.Range("H122").Value = 1708.091
$ws.Range("I122").Value = 1433.9
$ws.Range("J122").Value = 4450
$ws.Range("K122").Value = 4301.700000000001
$ws.Range("L122").Value = 13350
$ws.Range("M122").Value = -1851.700000000001
$ws.Range("N122").Value = -18250
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("H141").Value = 27500
$ws.Range("I141").Value = 13333.333
$ws.Range("K141").Value = 39999.999
$ws.Range("M141").Value = -34819.999
$ws.Range("M127").ClearContents()
$ws.Range("N127").ClearContents()

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 99921
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 99921
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 99921
$ws.Range("N133").Value = -104981
$ws.Range("H134").Value = 94987
$ws.Range("J134").Value = 94987
$ws.Range("L134").Value = 94987
$ws.Range("N134").Value = -105127
$ws.Range("M133").ClearContents()

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3312.25
$ws.Range("I134").Value = 3312.25
$ws.Range("K134").Value = 9936.75
$ws.Range("M134").Value = -7401.75

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 96.27273
$ws.Range("I40").Value = 69.28570999999999
$ws.Range("J40").Value = 143.5
$ws.Range("K40").Value = 277.14284
$ws.Range("L40").Value = 574
$ws.Range("M40").Value = -208.14284
$ws.Range("N40").Value = -712
$ws.Range("H120").Value = 1000
$ws.Range("I120").Value = 1000
$ws.Range("K120").Value = 3000
$ws.Range("M120").Value = 1838
$ws.Range("H132").Value = 2249.5
$ws.Range("J132").Value = 2249.5
$ws.Range("L132").Value = 20245.5
$ws.Range("N132").Value = -25305.5
$ws.Range("H140").Value = 1946.6666
$ws.Range("I140").Value = 1946.6666
$ws.Range("K140").Value = 5839.9998
$ws.Range("M140").Value = -659.9997999999996
$ws.Range("H141").Value = 17666
$ws.Range("I141").Value = 1999
$ws.Range("K141").Value = 5997
$ws.Range("M141").Value = -817

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 82.916664
$ws.Range("I2").Value = 98
$ws.Range("J2").Value = 61.8
$ws.Range("K2").Value = 98
$ws.Range("L2").Value = 61.8
$ws.Range("M2").Value = 15
$ws.Range("N2").Value = -287.8
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("H92").Value = 7990.2
$ws.Range("J92").Value = 8711.333000000001
$ws.Range("L92").Value = 8711.333000000001
$ws.Range("N92").Value = -12455.333
$ws.Range("H132").Value = 9198.6
$ws.Range("I132").Value = 8499.5
$ws.Range("J132").Value = 9664.666999999999
$ws.Range("K132").Value = 25498.5
$ws.Range("L132").Value = 28994.001
$ws.Range("M132").Value = -22968.5
$ws.Range("N132").Value = -34054.001
$ws.Range("N53").ClearContents()
$ws.Range("N63").ClearContents()
$ws.Range("N66").ClearContents()

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 10875
$ws.Range("J4").Value = 13333.333
$ws.Range("L4").Value = 13333.333
$ws.Range("N4").Value = -13559.333
$ws.Range("H28").Value = 10875
$ws.Range("J28").Value = 13333.333
$ws.Range("L28").Value = 13333.333
$ws.Range("N28").Value = -13797.333
$ws.Range("H37").Value = 10875
$ws.Range("J37").Value = 13333.333
$ws.Range("L37").Value = 13333.333
$ws.Range("N37").Value = -13547.333
$ws.Range("H40").Value = 7493.4346
$ws.Range("I40").Value = 8531.526
$ws.Range("J40").Value = 2562.5
$ws.Range("K40").Value = 8531.526
$ws.Range("L40").Value = 2562.5
$ws.Range("M40").Value = -8395.526
$ws.Range("N40").Value = -2834.5
$ws.Range("H61").Value = 4137.5
$ws.Range("I61").Value = 4183.3335
$ws.Range("K61").Value = 4183.3335
$ws.Range("M61").Value = -3981.3335
$ws.Range("H113").Value = 4137.5
$ws.Range("I113").Value = 4183.3335
$ws.Range("K113").Value = 4183.3335
$ws.Range("M113").Value = -2013.3335
$ws.Range("I122").Value = 5050
$ws.Range("K122").Value = 15150
$ws.Range("M122").Value = -12700
$ws.Range("H132").Value = 999
$ws.Range("I132").Value = 999
$ws.Range("K132").Value = 2997
$ws.Range("M132").Value = -467

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4030.5454
$ws.Range("I96").Value = 3333.8572
$ws.Range("J96").Value = 5249.75
$ws.Range("K96").Value = 3333.8572
$ws.Range("L96").Value = 5249.75
$ws.Range("M96").Value = -1960.8572
$ws.Range("N96").Value = -7995.75
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("H126").Value = 3342.7144
$ws.Range("I126").Value = 3650.6667
$ws.Range("K126").Value = 10952.0001
$ws.Range("M126").Value = -8482.000100000001
$ws.Range("N123").ClearContents()
